$d = $word.ActiveDocument

# The byline paragraph currently reads (across several runs):
#   "Nenad Vilendečić , Milana Grbić & Dragan Matić,"
# It needs a comma added right after "Grbić":
#   "Nenad Vilendečić , Milana Grbić, & Dragan Matić,"
#
# Locate the lone "ć" that ends "Grbić" (it sits in its own run, right
# before " & Dragan") and append a comma to just that run, leaving every
# other run untouched.
$cCedilla = [string][char]0x0107   # "ć"
$rng = $d.Content
$found = $rng.Find.Execute($cCedilla + " & Dragan", $true, $false, $false,
                            $false, $false, $true, 1, $false, $null, 0)
if ($found) {
    $insertPoint = $d.Range($rng.Start, $rng.Start + 1)
    $insertPoint.InsertAfter(",")
}

# Clean up the stray "_GoBack" bookmark that used to sit between the name
# and the ampersand.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
